$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Apply updated cryptocurrency price/volume data scraped on Fri Aug 23 00:54:35 UTC 2024

$ws.Range("D2").Value = '60.608.69'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = '2.628.47'
$ws.Range("E3").Value = '  +0.83%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.64'
$ws.Range("E5").Value = '  +2.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.59'
$ws.Range("E6").Value = '  +1.23%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.598'
$ws.Range("E8").Value = '  -0.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.53'
$ws.Range("E9").Value = '  +0.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.106'
$ws.Range("E10").Value = '  +0.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.377'
$ws.Range("E11").Value = '  +1.94%  '
$ws.Range("E12").Value = '  +1.04%  '
$ws.Range("D13").Value = '3.095.22'
$ws.Range("E13").Value = '  +0.76%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.74'
$ws.Range("E14").Value = '  +9.30%  '
$ws.Range("D15").Value = '60.618.09'
$ws.Range("E15").Value = '  -0.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000141'
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("D17").Value = '2.633.91'
$ws.Range("E17").Value = '  +0.45%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.48'
$ws.Range("E18").Value = '  +1.82%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.74'
$ws.Range("E19").Value = '  +1.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '349.41'
$ws.Range("E20").Value = '  +0.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.93'
$ws.Range("E21").Value = '  -1.95%  '
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("E23").Value = '  +0.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.72'
$ws.Range("E24").Value = '  -0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.994'
$ws.Range("E25").Value = '  -0.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.161'
$ws.Range("E26").Value = '  +1.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.08'
$ws.Range("E27").Value = '  +4.94%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.98'
$ws.Range("E28").Value = '  +8.31%  '
$ws.Range("D29").Value = '0.0₃0806'
$ws.Range("E29").Value = '  +1.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.57'
$ws.Range("E30").Value = '  +3.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '168.43'
$ws.Range("E31").Value = '  +4.76%  '
$ws.Range("E32").Value = '  +0.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.58'
$ws.Range("E33").Value = '  +0.28%  '
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.41'
$ws.Range("E34").Value = '  +4.05%  '
$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.02'
$ws.Range("E35").Value = '  +5.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.31'
$ws.Range("E36").Value = '  +6.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.64'
$ws.Range("E37").Value = '  +1.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '328.91'
$ws.Range("E38").Value = '  +11.17%  '
$ws.Range("B39").Value = 'Filecoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.04'
$ws.Range("E39").Value = '  +6.23%  '
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '38.64'
$ws.Range("E40").Value = '  +2.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.858'
$ws.Range("E41").Value = '  +0.80%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.27'
$ws.Range("E42").Value = '  +7.74%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '135.12'
$ws.Range("E43").Value = '  -3.81%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.28'
$ws.Range("E44").Value = '  +3.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0993'
$ws.Range("E45").Value = '  +1.23%  '
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  +0.36%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.612'
$ws.Range("E47").Value = '  +0.82%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0556'
$ws.Range("E48").Value = '  +1.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.16'
$ws.Range("E49").Value = '  +3.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0242'
$ws.Range("E50").Value = '  +0.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.76'
$ws.Range("E51").Value = '  +0.60%  '
